$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")

# Insert a new row above row 6 (existing rows 6.. shift down to 7..),
# carrying over formatting from the row above (matches D6/E6 s="2",
# AB6/AC6 s="5"/"6" in the target).
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new shipping-method test row.
$ws.Range("A6").Value = "StandardShipping method"
$ws.Range("R6").Value = "Standard (5 - 7 Business Days)"

# Switch the active tab from "My AccountPage" back to "DataSet", and move
# the selection/scroll position to R6.
$ws9 = $wb.Worksheets.Item("My AccountPage")
$ws9.Activate()
$ws.Activate()
$ws.Range("R6").Select()
